$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7141138105734797
$ws.Range("D2").Value = 0.1672121791461905
$ws.Range("E2").Value = 0.1622457546147498
$ws.Range("F2").Value = 1.551594903784846
$ws.Range("G2").Value = 0.00247424243004628
$ws.Range("I2").Value = 1.100511047659349
$ws.Range("J2").Value = 0.2012741135808693
$ws.Range("K2").Value = 0.4532032925836802
$ws.Range("L2").Value = 0.2558711733115899
$ws.Range("N2").Value = 1.959933792393025
$ws.Range("O2").Value = 3.822103738925733
$ws.Range("B3").Value = 0.6895132758629643
$ws.Range("D3").Value = 0.1667866191308747
$ws.Range("E3").Value = 0.1630340943234314
$ws.Range("F3").Value = 1.556305972820674
$ws.Range("G3").Value = 0.002476658153399547
$ws.Range("I3").Value = 1.109823531985771
$ws.Range("J3").Value = 0.2024220512801662
$ws.Range("K3").Value = 0.4040889223917929
$ws.Range("L3").Value = 0.2458726491909289
$ws.Range("N3").Value = 1.977695645413352
$ws.Range("O3").Value = 3.836778749004196
$ws.Range("B4").Value = 0.6746548597348863
$ws.Range("D4").Value = 0.1665749024031768
$ws.Range("E4").Value = 0.1635582576223937
$ws.Range("F4").Value = 1.559898821609622
$ws.Range("G4").Value = 0.002478222340150848
$ws.Range("I4").Value = 1.115961433256366
$ws.Range("J4").Value = 0.2031675815546481
$ws.Range("K4").Value = 0.3739438986220875
$ws.Range("L4").Value = 0.2398283726508765
$ws.Range("N4").Value = 1.989175535084975
$ws.Range("O4").Value = 3.847604189112175
$ws.Range("B5").Value = 0.6686625200555341
$ws.Range("D5").Value = 0.1665011473611813
$ws.Range("E5").Value = 0.1637819704806924
$ws.Range("F5").Value = 1.56153922373889
$ws.Range("G5").Value = 0.002478880167034349
$ws.Range("I5").Value = 1.118568372421223
$ws.Range("J5").Value = 0.2034816432977797
$ws.Range("K5").Value = 0.3616631692597707
$ws.Range("L5").Value = 0.2373893245016916
$ws.Range("N5").Value = 1.993998105901115
$ws.Range("O5").Value = 3.852472370686684
$ws.Range("B6").Value = 0.6676712949664818
$ws.Range("D6").Value = 0.1664896582398541
$ws.Range("E6").Value = 0.1638197292700321
$ws.Range("F6").Value = 1.561822264680046
$ws.Range("G6").Value = 0.002478990633002032
$ws.Range("I6").Value = 1.119007638806561
$ws.Range("J6").Value = 0.2035344127950516
$ws.Range("K6").Value = 0.3596242036632589
$ws.Range("L6").Value = 0.2369857802054725
$ws.Range("N6").Value = 1.994807615137423
$ws.Range("O6").Value = 3.853308324002015
$ws.Range("B7").Value = 0.674573790712941
$ws.Range("D7").Value = 0.1665738569465276
$ws.Range("E7").Value = 0.1635612337170507
$ws.Range("F7").Value = 1.559920230591267
$ws.Range("G7").Value = 0.002478231129164299
$ws.Range("I7").Value = 1.115996163286461
$ws.Range("J7").Value = 0.2031717755642899
$ws.Range("K7").Value = 0.3737782607950066
$ws.Range("L7").Value = 0.2397953811895235
$ws.Range("N7").Value = 1.989239989098945
$ws.Range("O7").Value = 3.847667993379901
$ws.Range("B8").Value = 0.7055807684708952
$ws.Range("D8").Value = 0.1670551880789759
$ws.Range("E8").Value = 0.1625092615417101
$ws.Range("F8").Value = 1.553074069288485
$ws.Range("G8").Value = 0.002475058611845799
$ws.Range("I8").Value = 1.103634871636778
$ws.Range("J8").Value = 0.2016614888726738
$ws.Range("K8").Value = 0.4362667249149297
$ws.Range("L8").Value = 0.2524041054859794
$ws.Range("N8").Value = 1.965938984106499
$ws.Range("O8").Value = 3.826787199146196
$ws.Range("B9").Value = 0.7683166388636096
$ws.Range("D9").Value = 0.1683902426571748
$ws.Range("E9").Value = 0.1607636346473988
$ws.Range("F9").Value = 1.545196755442959
$ws.Range("G9").Value = 0.002469476692692928
$ws.Range("I9").Value = 1.082723155880721
$ws.Range("J9").Value = 0.1990217796742582
$ws.Range("K9").Value = 0.5588699704672422
$ws.Range("L9").Value = 0.2778754576116711
$ws.Range("N9").Value = 1.924796536606717
$ws.Range("O9").Value = 3.800227848489328
$ws.Range("B10").Value = 0.8155590943359528
$ws.Range("D10").Value = 0.1696068272052997
$ws.Range("E10").Value = 0.159673142365989
$ws.Range("F10").Value = 1.542782263705888
$ws.Range("G10").Value = 0.002465761647447502
$ws.Range("I10").Value = 1.069383364057114
$ws.Range("J10").Value = 0.1972773669721408
$ws.Range("K10").Value = 0.6489579060724111
$ws.Range("L10").Value = 0.297036564272176
$ws.Range("N10").Value = 1.897337059149308
$ws.Range("O10").Value = 3.789472490426675
$ws.Range("B11").Value = 0.837295304634182
$ws.Range("D11").Value = 0.1702109012355777
$ws.Range("E11").Value = 0.1592184477253031
$ws.Range("F11").Value = 1.542414414670418
$ws.Range("G11").Value = 0.002464154585645487
$ws.Range("I11").Value = 1.063753145335145
$ws.Range("J11").Value = 0.1965258513612373
$ws.Range("K11").Value = 0.689938390837284
$ws.Range("L11").Value = 0.3058491531294436
$ws.Range("N11").Value = 1.885444699420097
$ws.Range("O11").Value = 3.786478717961529
$ws.Range("B12").Value = 0.8455609599535023
$ws.Range("D12").Value = 0.1704468791665406
$ws.Range("E12").Value = 0.1590521928381854
$ws.Range("F12").Value = 1.542379988163617
$ws.Range("G12").Value = 0.002463557897786811
$ws.Range("I12").Value = 1.061684050374609
$ws.Range("J12").Value = 0.1962472949674483
$ws.Range("K12").Value = 0.7054558449362673
$ws.Range("L12").Value = 0.3091999042566584
$ws.Range("N12").Value = 1.881027451407336
$ws.Range("O12").Value = 3.785617843252197
$ws.Range("B13").Value = 0.8437792737933023
$ws.Range("D13").Value = 0.1703957364835418
$ws.Range("E13").Value = 0.1590877354930491
$ws.Range("F13").Value = 1.542382741234974
$ws.Range("G13").Value = 0.002463685878144168
$ws.Range("I13").Value = 1.062126868286803
$ws.Range("J13").Value = 0.1963070193222691
$ws.Range("K13").Value = 0.7021139396906335
$ws.Range("L13").Value = 0.3084776577831008
$ws.Range("N13").Value = 1.881974954389639
$ws.Range("O13").Value = 3.785791119176508
$ws.Range("B14").Value = 0.8379746349379786
$ws.Range("D14").Value = 0.1702301707286793
$ws.Range("E14").Value = 0.1592046511515033
$ws.Range("F14").Value = 1.542409481647567
$ws.Range("G14").Value = 0.00246410525832199
$ws.Range("I14").Value = 1.063581658412314
$ws.Range("J14").Value = 0.1965028136815588
$ws.Range("K14").Value = 0.691215045225249
$ws.Range("L14").Value = 0.3061245497510896
$ws.Range("N14").Value = 1.885079563680135
$ws.Range("O14").Value = 3.786402427395899
$ws.Range("B15").Value = 0.8344236155081433
$ws.Range("D15").Value = 0.1701296966390871
$ws.Range("E15").Value = 0.1592770367406242
$ws.Range("F15").Value = 1.542439512818419
$ws.Range("G15").Value = 0.002464363685093605
$ws.Range("I15").Value = 1.064480955278992
$ws.Range("J15").Value = 0.196623527671052
$ws.Range("K15").Value = 0.6845390073040676
$ws.Range("L15").Value = 0.3046849701259049
$ws.Range("N15").Value = 1.886992442747772
$ws.Range("O15").Value = 3.78681239057147
$ws.Range("B16").Value = 0.8141434630302911
$ws.Range("D16").Value = 0.1695683639832808
$ws.Range("E16").Value = 0.1597036884610485
$ws.Range("F16").Value = 1.542820988193384
$ws.Range("G16").Value = 0.00246586833698618
$ws.Range("I16").Value = 1.069760123519188
$ws.Range("J16").Value = 0.1973273245005494
$ws.Range("K16").Value = 0.6462796448588506
$ws.Range("L16").Value = 0.2964625572448369
$ws.Range("N16").Value = 1.898126306489401
$ws.Range("O16").Value = 3.789706329166478
$ws.Range("B17").Value = 0.8017646236494897
$ws.Range("D17").Value = 0.1692369369677991
$ws.Range("E17").Value = 0.1599760073432606
$ws.Range("F17").Value = 1.543241987428672
$ws.Range("G17").Value = 0.00246681259394584
$ws.Range("I17").Value = 1.073110896342275
$ws.Range("J17").Value = 0.1977698328408155
$ws.Range("K17").Value = 0.6228079289953428
$ws.Range("L17").Value = 0.2914428424757176
$ws.Range("N17").Value = 1.905109999646289
$ws.Range("O17").Value = 3.791967836553539
$ws.Range("B18").Value = 0.7946677785169811
$ws.Range("D18").Value = 0.1690510789721955
$ws.Range("E18").Value = 0.1601365337360221
$ws.Range("F18").Value = 1.54355290522868
$ws.Range("G18").Value = 0.002467363515155065
$ws.Range("I18").Value = 1.075079413349961
$ws.Range("J18").Value = 0.1980283086806236
$ws.Range("K18").Value = 0.6093075770868097
$ws.Range("L18").Value = 0.2885646934574169
$ws.Range("N18").Value = 1.909183225756223
$ws.Range("O18").Value = 3.793447344027641
$ws.Range("B19").Value = 0.7922689024534577
$ws.Range("D19").Value = 0.1689889715664776
$ws.Range("E19").Value = 0.1601915549792778
$ws.Range("F19").Value = 1.543669993396655
$ws.Range("G19").Value = 0.002467551390159692
$ws.Range("I19").Value = 1.075753004248309
$ws.Range("J19").Value = 0.1981165042563253
$ws.Range("K19").Value = 0.6047366079100129
$ws.Range("L19").Value = 0.2875917636368683
$ws.Range("N19").Value = 1.910572035866215
$ws.Range("O19").Value = 3.793978987838756
$ws.Range("B20").Value = 0.8030799825482688
$ws.Range("D20").Value = 0.1692717246505495
$ws.Range("E20").Value = 0.1599466154929434
$ws.Range("F20").Value = 1.543190055169298
$ws.Range("G20").Value = 0.002466711268631644
$ws.Range("I20").Value = 1.072749932950327
$ws.Range("J20").Value = 0.1977223176888216
$ws.Range("K20").Value = 0.6253065435255678
$ws.Range("L20").Value = 0.291976263681093
$ws.Range("N20").Value = 1.904360737522261
$ws.Range("O20").Value = 3.791708597314425
$ws.Range("B21").Value = 0.8396786636531317
$ws.Range("D21").Value = 0.1702786056651036
$ws.Range("E21").Value = 0.1591701494738498
$ws.Range("F21").Value = 1.542398782614299
$ws.Range("G21").Value = 0.002463981754240467
$ws.Range("I21").Value = 1.06315264336456
$ws.Range("J21").Value = 0.1964451407170738
$ws.Range("K21").Value = 0.6944163475342009
$ws.Range("L21").Value = 0.3068153467475554
$ws.Range("N21").Value = 1.884165327251262
$ws.Range("O21").Value = 3.786215469634982
$ws.Range("B22").Value = 0.8637995271767522
$ws.Range("D22").Value = 0.1709787670007969
$ws.Range("E22").Value = 0.1586972275825467
$ws.Range("F22").Value = 1.54249282007072
$ws.Range("G22").Value = 0.002462267030510279
$ws.Range("I22").Value = 1.057247136272419
$ws.Range("J22").Value = 0.1956455495335083
$ws.Range("K22").Value = 0.7395775843133094
$ws.Range("L22").Value = 0.3165928086581715
$ws.Range("N22").Value = 1.871468487102716
$ws.Range("O22").Value = 3.78421537221854
$ws.Range("B23").Value = 0.8509075486201709
$ws.Range("D23").Value = 0.1706012428920474
$ws.Range("E23").Value = 0.1589464812988375
$ws.Range("F23").Value = 1.542386766157833
$ws.Range("G23").Value = 0.002463175899417302
$ws.Range("I23").Value = 1.060365465193481
$ws.Range("J23").Value = 0.1960690987313658
$ws.Range("K23").Value = 0.715475006043647
$ws.Range("L23").Value = 0.3113672083748327
$ws.Range("N23").Value = 1.878199099465329
$ws.Range("O23").Value = 3.785137463599881
$ws.Range("B24").Value = 0.8024852467084713
$ws.Range("D24").Value = 0.169255982546332
$ws.Range("E24").Value = 0.1599598911916313
$ws.Range("F24").Value = 1.543213319175194
$ws.Range("G24").Value = 0.00246675705271271
$ws.Range("I24").Value = 1.072912993316685
$ws.Range("J24").Value = 0.197743786604863
$ws.Range("K24").Value = 0.6241769389898764
$ws.Range("L24").Value = 0.2917350797928151
$ws.Range("N24").Value = 1.904699297612026
$ws.Range("O24").Value = 3.791825240731441
$ws.Range("B25").Value = 0.7511409538238922
$ws.Range("D25").Value = 0.1679874822131922
$ws.Range("E25").Value = 0.1612020509931416
$ws.Range("F25").Value = 1.546734781550029
$ws.Range("G25").Value = 0.002470918693519542
$ws.Range("I25").Value = 1.088024508475627
$ws.Range("J25").Value = 0.1997015572195791
$ws.Range("K25").Value = 0.5256985720699276
$ws.Range("L25").Value = 0.2709057037840665
$ws.Range("N25").Value = 1.935440140204122
$ws.Range("O25").Value = 3.805873927695444

Write-Output "applied 264 cells"
